$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.352.69"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "3.659.52"
$ws.Range("E3").Value = "  +1.79%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.76"
$ws.Range("E5").Value = "  -1.23%  "

$ws.Range("E6").Value = "  +14.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "660.55"
$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.424"
$ws.Range("E8").Value = "  +2.18%  "

$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").Value = "3.657.39"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.88"
$ws.Range("E12").Value = "  +2.27%  "

$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  +3.41%  "

$ws.Range("D15").Value = "4.339.42"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000271"
$ws.Range("E16").Value = "  +5.09%  "

$ws.Range("D17").Value = "96.194.75"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.90"
$ws.Range("E18").Value = "  +14.68%  "

$ws.Range("D19").Value = "3.650.97"
$ws.Range("E19").Value = "  +1.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.75"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.534"
$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "520.93"
$ws.Range("E23").Value = "  +1.96%  "

$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.89"
$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.16"
$ws.Range("E27").Value = "  +5.04%  "

$ws.Range("E28").Value = "  -1.74%  "

$ws.Range("E29").Value = "  +8.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.28"
$ws.Range("E30").Value = "  +5.96%  "

$ws.Range("E31").Value = "  -0.49%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +11.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.96"
$ws.Range("E35").Value = "  +4.20%  "

$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.587"
$ws.Range("E37").Value = "  +2.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "628.05"
$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "45.27"
$ws.Range("E39").Value = "  +35.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.71"
$ws.Range("E40").Value = "  -1.77%  "

$ws.Range("E41").Value = "  +4.79%  "

$ws.Range("E42").Value = "  +4.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.95"
$ws.Range("E43").Value = "  +5.19%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.27"
$ws.Range("E45").Value = "  +7.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0456"
$ws.Range("E46").Value = "  +5.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.436"
$ws.Range("E47").Value = "  +20.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.28"
$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.54"
$ws.Range("E50").Value = "  +2.46%  "

$ws.Range("B51").Value = "MantraDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.58"
$ws.Range("E51").Value = "  +1.64%  "
